$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated objective (B), gap (C), and solve time (D) values for rows 2-11
# per corrected fixed recourse data.
$values = @{
    2  = @(-123.47338606118794, 0.09136807702319806, 66.268335062)
    3  = @(-117.61674244158768, 0.0646728753907774,  86.711968143)
    4  = @(-122.00142694489112, 0.07604737257360063, 77.857640638)
    5  = @(-121.30725095787741, 0.09996849057912861, 98.849060709)
    6  = @(-120.12292887523154, 0.0840953498554457,  72.891467707)
    7  = @(-118.89272648754448, 0.08608142800931605, 65.708824145)
    8  = @(-115.03704432900317, 0.0,                 78.536159631)
    9  = @(-119.44147912757283, 0.09684089792055087, 102.95397379)
    10 = @(-120.14729823432504, 0.09943172612207225, 76.289173834)
    11 = @(-116.2650930538125,  0.0936922930803802,  84.065189864)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
